$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.002819337091903754
$ws.Range("B3").Value = 0.7518128636896441
$ws.Range("B4").Value = 0.09572526895753761
$ws.Range("B5").Value = "[0.9999999999999999, 0.01961579970753712, -0.04943807156094064, -0.08845748901224221, -0.06319791329113343, 0.08426028056372001, 0.041465496117445565, 0.06269606292350176, 0.02189001395998311, 0.03347874560399804, -0.036407113402379995, -0.04806140609213703, 0.00010435876312211298, 0.037290723862806144, 0.060547275675394856, 0.007070870783345106, 0.0037079404942592965, -0.035063526513363974, -0.01355569168909592, -0.012762033644481337]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1991060930081394
$ws.Range("B5").Value = 0.5266506103352919
$ws.Range("B6").Value = 0.04797501166961602
$ws.Range("B7").Value = 0.1362049765512335
$ws.Range("B8").Value = 0.6163612622924162
$ws.Range("B9").Value = 0.04668668767180976
$ws.Range("B10").Value = "[1.0, 0.09589449722006389, -0.00955350793759561, 0.0014675186852311718, 0.016972055997117447, 0.03566264772384684, -0.018150973590492846, 0.028603217436567668, -0.0037931927637253884, 0.02168541645559285, -0.0012983611908527354, 0.009175382740956447, 0.018189305728656808, 0.019514302607564957, 0.030378382594517892, 0.03140236853485812, 0.0007228280028319554, 0.02570631282528608, -0.012546011011400836, -0.022453455812434076]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = [double]"6.850956479803724e-05"
$ws.Range("B3").Value = 0.01557196050730895
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0.8570705910508745
$ws.Range("B6").Value = "[1.0, 0.04862741777659727, -0.04337281117969399, -0.11659102928073099, -0.0936167171415675, 0.11198320972385668, 0.06310276090488323, 0.07011427773845622, 0.054570794839879075, 0.024963632768077622, -0.059254500515258386, -0.057135089457150796, -0.0017462836607806775, 0.0672070867937096, 0.0660248408805556, 0.030487411083917403, 0.003337284532397537, -0.05395119855142194, -0.04799287651963639, -0.012161573988790544]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = [double]"7.774552107398647e-05"
$ws.Range("B3").Value = 0.01555262446802603
$ws.Range("B4").Value = [double]"3.414219926850568e-16"
$ws.Range("B5").Value = 0.8572793861801071
$ws.Range("B6").Value = "[1.0, 0.048650263111892564, -0.04335383026598802, -0.11657968737796189, -0.0936101431440052, 0.11198906383720478, 0.06311118157339579, 0.0701248729884461, 0.05458455327235358, 0.024970552615187908, -0.059245748704233464, -0.05712643280936157, -0.0017400869801134612, 0.06721146831027276, 0.06602645517773753, 0.03049113718010807, 0.0033414918901165007, -0.053943892796704807, -0.04798850000771013, -0.012155796885177284]"
$ws.Range("B7").Value = [double]"-6.895066249840756e-05"

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.002800769253995617
$ws.Range("B3").Value = 0.01340378975500062
$ws.Range("B4").Value = [double]"8.671616135718278e-10"
$ws.Range("B5").Value = 0.8595900300550963
$ws.Range("B6").Value = "[0.9999999999999999, 0.019195973676209004, -0.049400687144011565, -0.0885205651572413, -0.06355320353737376, 0.08380896211496351, 0.041455262804269664, 0.06245785920698182, 0.021986398175759356, 0.03369840271481277, -0.03643887109377647, -0.048151198595685556, -0.0003614505645174576, 0.03760413537291535, 0.060435621426453506, 0.007077469579872006, 0.0040685829305930074, -0.03482963149066882, -0.01357301682643516, -0.012974957284723272]"
$ws.Range("B7").Value = [double]"1.090133908609534e-10"
$ws.Range("B9").Value = 0.7520382037119598

